$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.615.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.81%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.563.39"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.65%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.20%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'210.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.88%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.486"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.59%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.17%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'24.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +5.24%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.244"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.61%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -0.47%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -0.07%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.788.48"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.51%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.560.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.98%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'28.645.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.93%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.514"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.54%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.52%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'61.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.33%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'231.28"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.46%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'7.36"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.77%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0674"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.49%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'  -0.19%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -1.23%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'8.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -0.57%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.13"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +3.99%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'150.66"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -0.34%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'14.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.97%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.37%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.16%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -2.37%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0460"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -4.48%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.87%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -0.93%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.388.54"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.24%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -4.81%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -2.96%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  -2.18%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.77%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.65"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.25%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.93%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +2.03%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -0.33%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.23%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.775"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -1.71%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -2.48%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'63.80"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.26%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'5.26"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.01%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.700.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.46%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -5.70%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'85.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.19%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'43.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.84%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  +0.40%  "
$ws.Range("E51").Style = "Normal"
